# Populate the "Retrived make" (E) and "Retrived colour" (F) columns by
# copying the values already present in "Known Make" (B) and
# "Known Colour" (C) for each data row (2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 7; $r++) {
    $make = $ws.Cells.Item($r, 2).Value2
    $colour = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 5).Value = $make
    $ws.Cells.Item($r, 6).Value = $colour
}
